$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.847.32'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.453.07'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '518.12'
$ws.Range('E5').Value = '  -3.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.38'
$ws.Range('E6').Value = '  -2.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.556'
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.458.77'
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0979'
$ws.Range('E10').Value = '  -3.31%  '
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.25'
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('E13').Value = '  -2.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.888.82'
$ws.Range('E14').Value = '  -1.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '57.788.40'
$ws.Range('E15').Value = '  -1.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.22'
$ws.Range('E16').Value = '  -3.19%  '
$ws.Range('E17').Value = '  -2.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.459.55'
$ws.Range('E18').Value = '  -2.45%  '
$ws.Range('E19').Value = '  -3.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '319.63'
$ws.Range('E20').Value = '  -1.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.16'
$ws.Range('E21').Value = '  -2.30%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.70'
$ws.Range('E23').Value = '  -4.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.36'
$ws.Range('E24').Value = '  -0.90%  '
$ws.Range('E25').Value = '  -2.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.994'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('E27').Value = '  -2.49%  '
$ws.Range('E28').Value = '  -2.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0740'
$ws.Range('E29').Value = '  -3.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.63'
$ws.Range('E30').Value = '  -1.79%  '
$ws.Range('E31').Value = '  -4.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.21'
$ws.Range('E32').Value = '  -6.41%  '
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('E37').Value = '  -5.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.95'
$ws.Range('E38').Value = '  -2.44%  '
$ws.Range('E39').Value = '  -3.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.18'
$ws.Range('E40').Value = '  -1.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.786'
$ws.Range('E41').Value = '  -2.66%  '
$ws.Range('E42').Value = '  -4.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '270.07'
$ws.Range('E43').Value = '  -4.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.95'
$ws.Range('E44').Value = '  -4.79%  '
$ws.Range('E45').Value = '  -3.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '124.56'
$ws.Range('E46').Value = '  -4.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0909'
$ws.Range('E47').Value = '  -1.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0486'
$ws.Range('E48').Value = '  -3.39%  '
$ws.Range('E49').Value = '  -3.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.76'
$ws.Range('E50').Value = '  -3.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.723.68'
$ws.Range('E51').Value = '  -1.71%  '
